$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: correct the Objetivos body text (was duplicated teacher name)
$ws.Range("B10").Value = "Apresentar aos alunos ingressantes o entendimento do que seja a carreira e as bases conceituais da Engenharia Física."
$ws.Range("C10").Value = "Apresentar aos alunos ingressantes o entendimento do que seja a carreira e as bases conceituais da Engenharia Física."

# Row 13: remove stray label, keep only first professor name in B/C
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"

# Row 14: remove stray label, keep only second professor name in B/C
$ws.Range("A14").Clear()
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# Row 15: "Programa resumido:" with its PT body text
$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = "A carreira de Engenharia Física. Conceitos básicos de Engenharia. Competências e habilidades de um engenheiro. Física conceitual. Realização de experimentos e projetos de Engenharia Física."
$ws.Range("C15").Value = "A carreira de Engenharia Física. Conceitos básicos de Engenharia. Competências e habilidades de um engenheiro. Física conceitual. Realização de experimentos e projetos de Engenharia Física."

# Row 16: "Short syllabus:" with its EN body text
$ws.Range("A16").Value = "Short syllabus:"
$ws.Range("B16").Value = "The Physics Engineering career. Basic engineering concepts. Skills and Abilities of an Engineer. Conceptual physics. Realization of experiments and projects of Physical Engineering."
$ws.Range("C16").Value = "The Physics Engineering career. Basic engineering concepts. Skills and Abilities of an Engineer. Conceptual physics. Realization of experiments and projects of Physical Engineering."

# Row 17: "Programa:" -- add B/C body text (new cells, copy formats first)
$ws.Range("A17").Value = "Programa:"
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("B17").Value = "A carreira de Engenharia Física. Cientistas x engenheiros: o papel interdisciplinar da Engenharia Física. Campos de atuação. A Física como ciência conceitual: Como aprender Física. Realização de demonstrações e experimentos científicos significativos de Física.Conceitos básicos de Engenharia. Habilidades e competências de um engenheiro.Desenvolvimento de um projeto temático de Engenharia Física.Competição entre projetos de diferentes grupos.Avaliação das competições e da disciplina como um todo."
$ws.Range("C17").Value = "A carreira de Engenharia Física. Cientistas x engenheiros: o papel interdisciplinar da Engenharia Física. Campos de atuação. A Física como ciência conceitual: Como aprender Física. Realização de demonstrações e experimentos científicos significativos de Física.Conceitos básicos de Engenharia. Habilidades e competências de um engenheiro.Desenvolvimento de um projeto temático de Engenharia Física.Competição entre projetos de diferentes grupos.Avaliação das competições e da disciplina como um todo."

# Row 18: "Syllabus:" with its EN body text
$ws.Range("A18").Value = "Syllabus:"
$ws.Range("B18").Value = "The career of Engineering Physics. Scientists x engineers: the interdisciplinary role of Engineering Physics. Fields of action.Physics as a conceptual science: How to learn Physics. Realization of demonstrations and significant scientific experiments in Physics.Basic engineering concepts. Skills and competences of an engineer.Development of a thematic project of Physical Engineering.Competition between projects from different groups.Evaluation of competitions and the discipline as a whole."
$ws.Range("C18").Value = "The career of Engineering Physics. Scientists x engineers: the interdisciplinary role of Engineering Physics. Fields of action.Physics as a conceptual science: How to learn Physics. Realization of demonstrations and significant scientific experiments in Physics.Basic engineering concepts. Skills and competences of an engineer.Development of a thematic project of Physical Engineering.Competition between projects from different groups.Evaluation of competitions and the discipline as a whole."

# Row 19: "Avaliacao:" -- A-only, clear stale B/C method text
$ws.Range("A19").Value = "Avaliação:"
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()

# Row 20: "Metodo:" with its body text
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "As atividades práticas e os projetos que serão desenvolvidos durante as aulas serão avaliados por docentes e pelos alunos (processo de avaliação crítica)."
$ws.Range("C20").Value = "As atividades práticas e os projetos que serão desenvolvidos durante as aulas serão avaliados por docentes e pelos alunos (processo de avaliação crítica)."

# Row 21: "Criterio:" with its body text
$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "A média final será uma composição de fatores relativos à participação do aluno nos trabalhos desenvolvidos, conjuntamente com o rendimento de seu grupo."
$ws.Range("C21").Value = "A média final será uma composição de fatores relativos à participação do aluno nos trabalhos desenvolvidos, conjuntamente com o rendimento de seu grupo."

# Row 22 (new): "Norma de recuperacao:" -- paste formats then set values
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "Devido às características da disciplina, não será oferecida recuperação."
$ws.Range("C22").Value = "Devido às características da disciplina, não será oferecida recuperação."

# Row 23 (new): "Bibliografia:" -- paste formats then set values
$ws.Range("A21").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "ARAÚJO-MOREIRA, F. M. Engenharia Física: a Carreira do Novo Milênio, São Carlos: Gráfica e Editora Guillen & Andriolli, 2014.`nBAZZO, A. B.; PEREIRA, L.T.V. Introdução à Engenharia. Editora da UFSC, Florianópolis, 1993.`nALEXANDER, C. K.; WATSON, J. A. Habilidades para uma carreira de sucesso na engenharia, Porto Alegre: AMGH Editora, 2015.`nBROCKMAN, J. B. Introdução à Engenharia. LTC, Rio de Janeiro, 2009.`nKNOWLEDGE FLOW. Engineering Physics - Ebook, Índia, 2015.`nCHAVES, A. S.; VALADARES, E. C.; ALVES, E. G. Aplicações da Física Quântica do Transistor à Nanotecnologia, São Paulo: Livraria da Física, 2005."
$ws.Range("C23").Value = "ARAÚJO-MOREIRA, F. M. Engenharia Física: a Carreira do Novo Milênio, São Carlos: Gráfica e Editora Guillen & Andriolli, 2014.`nBAZZO, A. B.; PEREIRA, L.T.V. Introdução à Engenharia. Editora da UFSC, Florianópolis, 1993.`nALEXANDER, C. K.; WATSON, J. A. Habilidades para uma carreira de sucesso na engenharia, Porto Alegre: AMGH Editora, 2015.`nBROCKMAN, J. B. Introdução à Engenharia. LTC, Rio de Janeiro, 2009.`nKNOWLEDGE FLOW. Engineering Physics - Ebook, Índia, 2015.`nCHAVES, A. S.; VALADARES, E. C.; ALVES, E. G. Aplicações da Física Quântica do Transistor à Nanotecnologia, São Paulo: Livraria da Física, 2005."

# Row heights for the final layout
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120

Write-Output "done"
